$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "67.110.97"
$ws.Range("E2").Value2 = "  +0.52%  "
$ws.Range("D3").Value2 = "3.875.49"
$ws.Range("E3").Value2 = "  +3.55%  "
$ws.Range("E4").Value2 = "  -0.14%  "
$ws.Range("D5").Value2 = "'428.10"
$ws.Range("E5").Value2 = "  +2.18%  "
$ws.Range("D6").Value2 = "'132.01"
$ws.Range("E6").Value2 = "  +0.91%  "
$ws.Range("D7").Value2 = "3.867.88"
$ws.Range("E7").Value2 = "  +3.61%  "
$ws.Range("D8").Value2 = "'0.615"
$ws.Range("E8").Value2 = "  -5.71%  "
$ws.Range("E9").Value2 = "  -0.09%  "
$ws.Range("D10").Value2 = "'0.733"
$ws.Range("E10").Value2 = "  -4.73%  "
$ws.Range("E11").Value2 = "  -6.80%  "
$ws.Range("D12").Value2 = "'0.0000365"
$ws.Range("E12").Value2 = "  -9.04%  "
$ws.Range("E13").Value2 = "  -4.17%  "
$ws.Range("D14").Value2 = "4.491.35"
$ws.Range("E14").Value2 = "  +3.96%  "
$ws.Range("D15").Value2 = "'10.13"
$ws.Range("E15").Value2 = "  -3.42%  "
$ws.Range("D16").Value2 = "'15.60"
$ws.Range("E16").Value2 = "  +18.07%  "
$ws.Range("D17").Value2 = "3.895.00"
$ws.Range("E17").Value2 = "  +4.00%  "
$ws.Range("E18").Value2 = "  -1.09%  "
$ws.Range("D19").Value2 = "'19.64"
$ws.Range("E19").Value2 = "  -5.28%  "
$ws.Range("D20").Value2 = "67.244.13"
$ws.Range("E20").Value2 = "  +0.58%  "
$ws.Range("E21").Value2 = "  -5.94%  "
$ws.Range("D22").Value2 = "'408.65"
$ws.Range("E22").Value2 = "  -8.06%  "
$ws.Range("D23").Value2 = "'14.55"
$ws.Range("E23").Value2 = "  -12.14%  "
$ws.Range("D24").Value2 = "'85.31"
$ws.Range("E24").Value2 = "  -4.88%  "
$ws.Range("D26").Value2 = "'37.63"
$ws.Range("E26").Value2 = "  -2.74%  "
$ws.Range("D27").Value2 = "'5.69"
$ws.Range("E27").Value2 = "  +11.58%  "
$ws.Range("D28").Value2 = "'3.22"
$ws.Range("E28").Value2 = "  -3.26%  "
$ws.Range("E29").Value2 = "  -5.88%  "
$ws.Range("D30").Value2 = "'691.29"
$ws.Range("E30").Value2 = "  +4.63%  "
$ws.Range("E31").Value2 = "  -1.53%  "
$ws.Range("D32").Value2 = "'12.47"
$ws.Range("E32").Value2 = "  -1.83%  "
$ws.Range("D33").Value2 = "'2.76"
$ws.Range("E33").Value2 = "  +0.12%  "
$ws.Range("D34").Value2 = "'7.15"
$ws.Range("E34").Value2 = "  -1.16%  "
$ws.Range("E35").Value2 = "  -7.25%  "
$ws.Range("E36").Value2 = "  -8.03%  "
$ws.Range("D37").Value2 = "0.0₃0812"
$ws.Range("E37").Value2 = "  +8.82%  "
$ws.Range("E38").Value2 = "  -0.06%  "
$ws.Range("D39").Value2 = "'55.35"
$ws.Range("E39").Value2 = "  -3.01%  "
$ws.Range("B40").Value2 = "VeChain"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value2 = "'0.0459"
$ws.Range("E40").Value2 = "  -6.57%  "
$ws.Range("B41").Value2 = "ThetaToken"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value2 = "'3.05"
$ws.Range("E41").Value2 = "  +0.96%  "
$ws.Range("E43").Value2 = "  -9.38%  "
$ws.Range("B44").Value2 = "NEARProtocol"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value2 = "'4.50"
$ws.Range("E44").Value2 = "  +3.74%  "
$ws.Range("B45").Value2 = "Monero"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value2 = "'147.74"
$ws.Range("E45").Value2 = "  +0.92%  "
$ws.Range("B46").Value2 = "EnergySwap"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value2 = "'26.57"
$ws.Range("E46").Value2 = "  -8.88%  "
$ws.Range("E47").Value2 = "  -3.02%  "
$ws.Range("D48").Value2 = "'3.26"
$ws.Range("E48").Value2 = "  -5.58%  "
$ws.Range("D49").Value2 = "'3.11"
$ws.Range("E49").Value2 = "  -4.89%  "
$ws.Range("D50").Value2 = "'2.79"
$ws.Range("E50").Value2 = "  -3.45%  "
$ws.Range("E51").Value2 = "  -4.98%  "
